$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 89
$prev = 88

# --- Values ---
$ws.Cells.Item($row, 1).Value = 88
$ws.Cells.Item($row, 2).Value = "denmark"
$ws.Cells.Item($row, 3).Value = "1st-division"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45235.625
$ws.Cells.Item($row, 6).Value = "Aalborg"
$ws.Cells.Item($row, 7).Value = 2
$ws.Cells.Item($row, 8).Value = "Helsingor"
$ws.Cells.Item($row, 9).Value = 2
$ws.Cells.Item($row, 10).Value = 1.22
$ws.Cells.Item($row, 11).Value = "29/10/2023 15:12"
$ws.Cells.Item($row, 12).Value = 1.23
$ws.Cells.Item($row, 13).Value = "05/11/2023 14:04"
$ws.Cells.Item($row, 14).Value = 7.27
$ws.Cells.Item($row, 15).Value = "29/10/2023 15:12"
$ws.Cells.Item($row, 16).Value = 6.59
$ws.Cells.Item($row, 17).Value = "05/11/2023 14:55"
$ws.Cells.Item($row, 18).Value = 10.83
$ws.Cells.Item($row, 19).Value = "29/10/2023 15:12"
$ws.Cells.Item($row, 20).Value = 11.31
$ws.Cells.Item($row, 21).Value = "05/11/2023 14:55"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/denmark/1st-division/aalborg-helsingor/rJCDn5ra/"

# --- Formatting: replicate the styled columns (A = bold/border/centered, E = datetime) ---
# using copy/paste-special (formats only) from the row above so the existing
# style entries in styles.xml are reused instead of new ones being minted.
$ws.Cells.Item($prev, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($prev, 5).Copy() | Out-Null
$ws.Cells.Item($row, 5).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
